# Update the 苏州-漫展信息 workbook with refreshed "想去人数" (interested count)
# and "最低票价" (min ticket price) figures on the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 114
$ws1.Range("F6").Value = 71
$ws1.Range("F9").Value = 1281
$ws1.Range("G10").Value = 58
$ws1.Range("F11").Value = 1007
$ws1.Range("F12").Value = 10374
$ws1.Range("F13").Value = 4
$ws1.Range("F17").Value = 661
$ws1.Range("F18").Value = 11908
$ws1.Range("F19").Value = 12299

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 114
$ws4.Range("F6").Value = 71
$ws4.Range("F10").Value = 1281
$ws4.Range("G11").Value = 58
$ws4.Range("F12").Value = 1007
$ws4.Range("F13").Value = 10374
$ws4.Range("F14").Value = 4
$ws4.Range("F18").Value = 661
$ws4.Range("F19").Value = 11908
$ws4.Range("F20").Value = 12299
